$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FT fuel - Diesel")

# Remove the "market group for electricity, low voltage" exchange row
# (row 229) from the hydrogen-production activity's exchange list - this
# was causing energy efficiency values above 1 for this activity.
$ws.Rows.Item(229).Delete()

# Reflect where the editor ended up after making the change.
$ws.Range("B232").Select()
